# LOQ4238.docx edit script
# Applies the content re-shuffle described by the diff:
#  - Objetivos (EN) paragraph gets the "Identify/Design" text (was "Develop...").
#  - Docente(s) paragraph gets the "Identificar/Projetar" text moved to the front,
#    loses the "11079086 - Herlandi..." run, and gains the "A recuperacao e
#    continua..." sentence appended at the end.
#  - "Programa resumido" paragraph becomes the "Livros e Artigos..." sentence.
#  - The second EN italic paragraph becomes the "Develop an inter/transdisciplinary..." text.
#  - "Programa" paragraph becomes the "O grupo social alvo..." sentence.
#  - Avaliacao paragraph loses "O grupo social alvo.../Criterio:" right after
#    "Metodo:", renames the old "Norma de recuperacao:" run to "Criterio:", and
#    appends a new bold "Norma de recuperacao:" run + the satisfaction-survey text.
#  - "Bibliografia" paragraph becomes "11079086 - Herlandi de Souza Andrade".

$d = $word.ActiveDocument
$BR = [char]11

function Find-Replace($range, [string]$old, [string]$new) {
    return $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

function Append-Run($paragraph, [string]$text) {
    $r = $paragraph.Range.Duplicate
    $r.Collapse(0)          # wdCollapseEnd
    $r.MoveEnd(1, -1)        # exclude the paragraph mark
    $r.Collapse(0)
    $insStart = $r.End
    $r.InsertAfter($text)
    $insEnd = $insStart + $text.Length
    return $d.Range($insStart, $insEnd)
}

# ----------------------------------------------------------------------------
# 1) Objetivos (EN, italic) paragraph -> "Identify a real...Design a production system."
# ----------------------------------------------------------------------------
$pObjEn = $d.Paragraphs.Item(7)
$oldObjEn = "- Develop an inter/transdisciplinary project, of medium-high complexity, on a topic related to Production Engineering, similar to situations that students will encounter in real life, in the actual exercise of their profession; Apply and integrate knowledge acquired in other course subjects; Develop technical skills (related to the project itself), and transversal skills (active learning, systemic thinking, problem-solving skills, resilience, analytical skills, teamwork, leadership, commitment, interpersonal relationships, conflict management, communication skills, planning skills, creativity, innovation, and initiative), in a learning environment based on PBL (Project-Based Learning and Problem -Baed Learning)."
$newObjEn = "- Identify a real or potential problem in a production system and propose a solution to the problem." + $BR + $BR + "- Design a production system."
Find-Replace $pObjEn.Range $oldObjEn $newObjEn | Out-Null

# ----------------------------------------------------------------------------
# 2) Docente(s) Responsavel(eis) paragraph
# ----------------------------------------------------------------------------
$pDoc = $d.Paragraphs.Item(9)

# 2a) remove "11079086 - Herlandi de Souza Andrade" + its trailing break entirely
$oldDocA = "11079086 - Herlandí de Souza Andrade" + $BR
Find-Replace $pDoc.Range $oldDocA "" | Out-Null

# 2b) remove the "Identificar.../Projetar..." run from its old spot (it will be
#     reinserted at the front of the paragraph)
$movedText = "- Identificar um problema real ou potencial em um sistema produtivo e propor uma solução para o problema." + $BR + $BR + "- Projetar um sistema produtivo." + $BR
$pDoc = $d.Paragraphs.Item(9)
Find-Replace $pDoc.Range $movedText "" | Out-Null

# 2c) insert that text back in as a brand new run at the very start of the paragraph
$pDoc = $d.Paragraphs.Item(9)
$startRange = $pDoc.Range.Duplicate
$startRange.Collapse(1)   # wdCollapseStart
$startRange.InsertBefore($movedText)

# 2d) append a break + new run with the "A recuperacao e continua..." sentence
$pDoc = $d.Paragraphs.Item(9)
Find-Replace $pDoc.Range "dentre outros." ("dentre outros." + $BR) | Out-Null
$pDoc = $d.Paragraphs.Item(9)
Append-Run $pDoc "A recuperação é contínua ao longo da disciplina, considerando as diversas atividades e entregas a serem realizadas. Não há prova de recuperação." | Out-Null

# ----------------------------------------------------------------------------
# 3) "Programa resumido" paragraph -> "Livros e Artigos cientificos..."
# ----------------------------------------------------------------------------
$pResumo = $d.Paragraphs.Item(11)
Find-Replace $pResumo.Range "A recuperação é contínua ao longo da disciplina, considerando as diversas atividades e entregas a serem realizadas. Não há prova de recuperação." "Livros e Artigos científicos relacionados com o tema do projeto/problema." | Out-Null

# ----------------------------------------------------------------------------
# 4) Second EN italic paragraph -> "Develop an inter/transdisciplinary project..."
# ----------------------------------------------------------------------------
$pObjEn2 = $d.Paragraphs.Item(12)
$oldObjEn2 = "- Identify a real or potential problem in a production system and propose a solution to the problem." + $BR + $BR + "- Design a production system."
$newObjEn2 = "- Develop an inter/transdisciplinary project, of medium-high complexity, on a topic related to Production Engineering, similar to situations that students will encounter in real life, in the actual exercise of their profession; Apply and integrate knowledge acquired in other course subjects; Develop technical skills (related to the project itself), and transversal skills (active learning, systemic thinking, problem-solving skills, resilience, analytical skills, teamwork, leadership, commitment, interpersonal relationships, conflict management, communication skills, planning skills, creativity, innovation, and initiative), in a learning environment based on PBL (Project-Based Learning and Problem -Baed Learning)."
Find-Replace $pObjEn2.Range $oldObjEn2 $newObjEn2 | Out-Null

# ----------------------------------------------------------------------------
# 5) "Programa" paragraph -> "O grupo social alvo sao medias..."
# ----------------------------------------------------------------------------
$pPrograma = $d.Paragraphs.Item(14)
Find-Replace $pPrograma.Range "Livros e Artigos científicos relacionados com o tema do projeto/problema." "O grupo social alvo são médias e grandes empresas, incluindo os profissionais dessas empresas, da Região do Vale do Paraíba." | Out-Null

# ----------------------------------------------------------------------------
# 6) Avaliacao paragraph
# ----------------------------------------------------------------------------
$pAval = $d.Paragraphs.Item(17)

# 6a) remove the "O grupo social alvo.../Paraiba." run + break (now located right
#     after "Metodo: ")
$oldAvalA = "O grupo social alvo são médias e grandes empresas, incluindo os profissionais dessas empresas, da Região do Vale do Paraíba." + $BR
Find-Replace $pAval.Range $oldAvalA "" | Out-Null

# 6b) remove the bold "Criterio: " run that followed it
$pAval = $d.Paragraphs.Item(17)
Find-Replace $pAval.Range "Critério: " "" | Out-Null

# 6c) rename the bold "Norma de recuperacao: " run (that precedes "A atividade
#     consiste...") to "Criterio: "
$pAval = $d.Paragraphs.Item(17)
Find-Replace $pAval.Range "Norma de recuperação: " "Critério: " | Out-Null

# 6d) append a break, then a new bold "Norma de recuperacao: " run, then the
#     satisfaction-survey sentence as a plain run
$pAval = $d.Paragraphs.Item(17)
Find-Replace $pAval.Range "lições aprendidas." ("lições aprendidas." + $BR) | Out-Null

$pAval = $d.Paragraphs.Item(17)
$normaText = "Norma de recuperação: "
$normaRange = Append-Run $pAval $normaText
$normaRange.Bold = 1

$pAval = $d.Paragraphs.Item(17)
$pesquisaText = "Será realizada uma pesquisa de satisfação com os responsáveis pela empresa participante da atividade, durante e após o projeto. Após a pesquisa, o grupo de estudantes da disciplina, fará uma análise dos resultados e uma autoavaliação e discutirá tais resultados com o professor da disciplina, de maneira e retroalimentar a realização do projeto."
Append-Run $pAval $pesquisaText | Out-Null

# ----------------------------------------------------------------------------
# 7) "Bibliografia" paragraph -> "11079086 - Herlandi de Souza Andrade"
# ----------------------------------------------------------------------------
$pBib = $d.Paragraphs.Item(19)
Find-Replace $pBib.Range "Será realizada uma pesquisa de satisfação com os responsáveis pela empresa participante da atividade, durante e após o projeto. Após a pesquisa, o grupo de estudantes da disciplina, fará uma análise dos resultados e uma autoavaliação e discutirá tais resultados com o professor da disciplina, de maneira e retroalimentar a realização do projeto." "11079086 - Herlandí de Souza Andrade" | Out-Null

Write-Output "All edits applied."
